$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the speed-conversion table (D7:J8) ---
$ws.Range("D8").Value = 75
$ws.Range("H8").Value = 1

# --- Remove the old "Conversão" / "Distancia" block that lived at H12:H13 ---
$ws.Range("H12:J12").UnMerge()
$ws.Range("H12:J12").ClearContents()
$ws.Range("H13").ClearContents()

# --- Rename the now-unused shared string so it can be reused as "Anda" ---
# (done implicitly below by typing "Anda" into F11)

# --- Build the new Distancia / % / Anda table ---
$ws.Range("D11").Value = "Distancia"
$ws.Range("E11").Value = "%"
$ws.Range("F11").Value = "Anda"

$ws.Range("D12").Value = 50
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 20

$ws.Range("D13").Value = 40
$ws.Range("E13").Value = 55
$ws.Range("F13").Value = 22

$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 60
$ws.Range("F14").Value = 24

$ws.Range("D15").Value = 20
$ws.Range("E15").Value = 65
$ws.Range("F15").Value = 26

$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 70
$ws.Range("F16").Value = 28

$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 30

# Center-align the table (mirrors the existing table's formatting)
$ws.Range("D11:G12").HorizontalAlignment = -4108
$ws.Range("D13:F20").HorizontalAlignment = -4108

# --- Helper interpolation values/formulas off to the side ---
$ws.Range("K11").Value = 50
$ws.Range("L11").Value = 50
$ws.Range("K12").Value = 25
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 75
$ws.Range("K16").Formula = "=((K12-K11)/(K13-K11))*(L13-L11)+L11"
$ws.Range("N16").Formula = "=(K12-K11)/(K13-K11)"

# --- Leave the selection where the last edit happened ---
$ws.Range("N16").Select() | Out-Null
